# Applies the "homeworks from 4 students" edit to the listening/attendance
# roster workbook.
#
# Summary of the change:
#  - The listening date in H1 was corrected (2015-02-22 -> 2015-03-22).
#  - A new listening date column (N1, 2015-05-17) plus a homework /
#    homework-grade pair of columns (O1/P1) were added.
#  - Several students got a new "n" (absent) mark for the new date columns.
#  - A stray duplicate mark in K9 was removed.
#  - A note column (N) got a few "n"/"?" remarks, and M15's old "qingjia"
#    note was normalized to the Chinese "请假".
#  - 4 students (rows 24, 36, 56, 59) received homework grades in column O.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: fix the first listening date, add new date/homework cols ---
$ws.Range("H1").Value = [DateTime]::new(2015, 3, 22)
$ws.Range("N1").Value = [DateTime]::new(2015, 5, 17)
$ws.Range("N1").NumberFormat = $ws.Range("M1").NumberFormat
$ws.Range("O1").Value = "作业"
$ws.Range("P1").Value = "作业成绩"

# --- Row 9 (范凌云): drop the stray K9 mark, add a leave note in N9 ---
$ws.Range("K9").ClearContents()
$ws.Range("N9").Value = "请假"

# --- New-date absence marks ("n") for students who missed 2015-05-17 ---
$ws.Range("H12").Value = "n"
$ws.Range("H15").Value = "n"
$ws.Range("H17").Value = "n"
$ws.Range("H18").Value = "n"
$ws.Range("H19").Value = "n"
$ws.Range("H27").Value = "n"
$ws.Range("H28").Value = "n"
$ws.Range("H30").Value = "n"
$ws.Range("H32").Value = "n"
$ws.Range("H40").Value = "n"
$ws.Range("H42").Value = "n"
$ws.Range("H49").Value = "n"
$ws.Range("H65").Value = "n"

# --- Note column (N) remarks ---
$ws.Range("N12").Value = "n"
$ws.Range("N19").Value = "?"
$ws.Range("N36").Value = "n"
$ws.Range("N47").Value = "?"
$ws.Range("N58").Value = "n"

# --- Normalize the old pinyin "qingjia" note to Chinese ---
$ws.Range("M15").Value = "请假"

# --- Homework grades for the 4 students who turned theirs in ---
$ws.Range("O59").Value = 1
$ws.Range("O24").Value = 2
$ws.Range("O56").Value = 3
$ws.Range("O36").Value = 4
